$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI output values ("Natmi following Dr Hou advice")
# Row 2: Sending cluster = ECs
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 261.380203
$ws.Range("H2").Value = 784.1406089999999
$ws.Range("I2").Value = 0.6968677182772199
$ws.Range("J2").Value = 0.6968677182772199
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.321929333333333
$ws.Range("N2").Value = 21.965788
$ws.Range("Q2").Value = 1913.807375498321
$ws.Range("R2").Value = 17224.26637948489
$ws.Range("S2").Value = 0.6968677182772199
$ws.Range("T2").Value = 0.6968677182772199

# Row 3: Sending cluster = FAPs
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 31.999428
$ws.Range("H3").Value = 95.998284
$ws.Range("I3").Value = 0.08531391482826334
$ws.Range("J3").Value = 0.08531391482826335
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.321929333333333
$ws.Range("N3").Value = 21.965788
$ws.Range("Q3").Value = 234.297550523088
$ws.Range("R3").Value = 2108.677954707792
$ws.Range("S3").Value = 0.08531391482826334
$ws.Range("T3").Value = 0.08531391482826335

# Row 4: Sending cluster = sCs
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 81.69901900000001
$ws.Range("H4").Value = 245.097057
$ws.Range("I4").Value = 0.2178183668945166
$ws.Range("J4").Value = 0.2178183668945167
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.321929333333333
$ws.Range("N4").Value = 21.965788
$ws.Range("Q4").Value = 598.1944437206574
$ws.Range("R4").Value = 5383.749993485916
$ws.Range("S4").Value = 0.2178183668945166
$ws.Range("T4").Value = 0.2178183668945167
